$d = $word.ActiveDocument

# --- Locate the last paragraph in the document (the one that currently
#     ends with "... iron out the specs for Gr. 8 night" and carries the
#     _GoBack bookmark). ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n - 1)

# Remove the existing _GoBack bookmark; we will recreate it at the end of
# the new final bullet once the new bullets have been appended.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Append the three new bullet paragraphs. InsertParagraphAfter()
#     inherits the pPr (ListParagraph style + numPr) of the paragraph it
#     is called on, which is exactly the formatting the other bullets in
#     this list use. ---
$newTexts = @(
    "I have wasted too much time trying to synchronize info when players are off their screen",
    "I am going to leave this feature for later and come back to it if I have the time",
    "Will move on to creating a mobile controller for the game"
)

$curPara = $lastPara
foreach ($t in $newTexts) {
    $curPara.Range.InsertParagraphAfter()
    $n = $d.Paragraphs.Count
    $curPara = $d.Paragraphs.Item($n - 1)
    $curPara.Range.Text = $t
}

# --- Re-create the _GoBack bookmark at the very end of the text of the
#     last bullet (right after "...the game", before the paragraph
#     mark). A zero-length Range positioned exactly at
#     "paragraph end - 1" cannot be used directly to seed
#     Bookmarks.Add, so a throw-away marker character is inserted,
#     wrapped in the bookmark, and then removed again - this leaves the
#     bookmark collapsed in the correct spot. ---
$endPos = $curPara.Range.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertAfter("X")

$endPos2 = $curPara.Range.End
$wrapRange = $d.Range($endPos2 - 2, $endPos2 - 1)
$d.Bookmarks.Add("_GoBack", $wrapRange)

$delRange = $d.Range($endPos2 - 2, $endPos2 - 1)
$delRange.Delete()
